$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "removed ICPetcdHighNumberOfFailedGRPCRequests because of false alerts on ICP"
# That alert has two rows in the sheet (warning + critical severity). Locate
# them dynamically by Alert Name (column A) and delete both entire rows.
$alertName = "ICPetcdHighNumberOfFailedGRPCRequests"
$alertCol = $ws.Columns.Item(1)

$rowsToDelete = @()
$first = $alertCol.Find($alertName)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $cell = $first
    do {
        $rowsToDelete += $cell.Row
        $cell = $alertCol.FindNext($cell)
    } while ($cell -ne $null -and $cell.Address() -ne $firstAddress)
}

# Delete from the bottom up so row numbers found above stay valid.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

# Reflect the author's final selection state: the row that now sits where the
# first deleted row used to be (now holding ICPetcdHighNumberOfFailedProposals).
if ($rowsToDelete.Count -gt 0) {
    $selectRow = ($rowsToDelete | Measure-Object -Minimum).Minimum
    $ws.Rows.Item($selectRow).Select() | Out-Null
}
